# "Fin 1ere version mapping" -- add the OncoFAIR MR Element Rank extension
# mapping column to the Elements table, and refresh the generation
# timestamp on the Metadata sheet.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: the IG-export timestamp moved forward ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value2 = "2024-04-22T13:59:04+00:00"

# --- Elements sheet: insert a new mapping column just before the
#     existing "Mapping: RIM Mapping" column (column AK), pushing the
#     RIM mapping column to AL ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Columns("AK:AK").Insert()

$elements.Range("AK1").Value2 = "Mapping: Mapping de l'extension OncoFAIR MR Element Rank"
$elements.Range("AK2").Value2 = "ELEMENT DE PRESCRIPTION"
$elements.Range("AK6").Value2 = "Rang élément prescription"

# Resize the new column to fit its (longer) contents, like Excel's
# "best fit" would when the sheet was last generated.
$elements.Columns("AK:AK").ColumnWidth = 64.2
